# Update ticket/sales counts (column F) on the "展览" and "全部类型" sheets.
# These figures were refreshed by the site's data generator (gh-pages rebuild),
# so only the raw numeric values in column F change; everything else is untouched.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1.xml) updates: row => new F value
$sheet1Updates = @{
    3  = 3139
    5  = 2153
    8  = 920
    9  = 996
    10 = 237
    11 = 457
    14 = 75
    16 = 7694
    17 = 335
    22 = 459
    23 = 529
    26 = 976
    28 = 1647
    30 = 1164
    32 = 479
    36 = 37
    37 = 161
    38 = 330
}

foreach ($row in $sheet1Updates.Keys) {
    $sheet1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheet4.xml) updates: row => new F value
$sheet4Updates = @{
    5  = 3139
    7  = 2153
    10 = 920
    12 = 996
    13 = 237
    14 = 457
    17 = 75
    19 = 7695
    20 = 335
    26 = 459
    27 = 529
    30 = 976
    32 = 1647
    34 = 1164
    36 = 479
    40 = 37
    41 = 161
    42 = 330
}

foreach ($row in $sheet4Updates.Keys) {
    $sheet4.Range("F$row").Value = $sheet4Updates[$row]
}
